# Auto-generated edit script applying the diff to before.xlsx
# Updates ticket-sales counts (col F / col G) and a block of shifted
# exhibition rows (a new con was inserted, cascading rows 38-48 down)
# on the '展览' sheet, plus numeric tweaks across the other sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 1889
$ws.Range('F5').Value = 765
$ws.Range('F6').Value = 120
$ws.Range('F8').Value = 928
$ws.Range('F9').Value = 1626
$ws.Range('F10').Value = 1288
$ws.Range('F11').Value = 1559
$ws.Range('F12').Value = 71
$ws.Range('F13').Value = 1557
$ws.Range('F14').Value = 348
$ws.Range('F15').Value = 1701
$ws.Range('F16').Value = 813
$ws.Range('F17').Value = 1132
$ws.Range('F18').Value = 378
$ws.Range('F19').Value = 58
$ws.Range('F21').Value = 1813
$ws.Range('F22').Value = 252
$ws.Range('F23').Value = 825
$ws.Range('F24').Value = 8
$ws.Range('F26').Value = 1260
$ws.Range('F27').Value = 1077
$ws.Range('F28').Value = 82
$ws.Range('F29').Value = 579
$ws.Range('F30').Value = 1175
$ws.Range('F31').Value = 912
$ws.Range('F33').Value = 1178
$ws.Range('F34').Value = 1120
$ws.Range('F35').Value = 286
$ws.Range('F36').Value = 86
$ws.Range('G36').Value = 68
$ws.Range('F37').Value = 891
$ws.Range('C38').Value = '上海·星芒旋转 Anikura动漫展'
$ws.Range('D38').Value = '海潮路133号B1 JUMP工坊'
$ws.Range('E38').Value = '2024.07.20 14:00-07.20 19:00'
$ws.Range('F38').Value = 3
$ws.Range('H38').Value = 'https://show.bilibili.com/platform/detail.html?id=87684'
$ws.Range('I38').Value = '//i1.hdslb.com/bfs/openplatform/202406/jzbMyROi1718612232328.jpeg'
$ws.Range('B39').Value = '2024-07-20'
$ws.Range('C39').Value = '上海·第九届Redamancy动漫游戏嘉年华'
$ws.Range('D39').Value = '中山北路3300号4楼 上海环球港'
$ws.Range('E39').Value = '2024.07.20 10:00-07.21 17:00'
$ws.Range('F39').Value = 1699
$ws.Range('H39').Value = 'https://show.bilibili.com/platform/detail.html?id=84637'
$ws.Range('I39').Value = '//i1.hdslb.com/bfs/openplatform/202404/hWLkXqwM1713194236349.png'
$ws.Range('B40').Value = '2024-07-27'
$ws.Range('C40').Value = '上海·第十二届Redamancy动漫游戏嘉年华'
$ws.Range('D40').Value = '中山北路3300号4楼L4001号 环球港上海世嘉都市乐园'
$ws.Range('E40').Value = '2024.07.27 10:00-07.28 17:00'
$ws.Range('F40').Value = 12
$ws.Range('H40').Value = 'https://show.bilibili.com/platform/detail.html?id=87423'
$ws.Range('I40').Value = '//i0.hdslb.com/bfs/openplatform/202406/Ll5fKZSD1718292388904.png'
$ws.Range('C41').Value = '上海·LOVELIVE ONLY'
$ws.Range('D41').Value = '海潮路133号B1 JUMP工坊'
$ws.Range('E41').Value = '2024.08.03 14:00-08.03 19:00'
$ws.Range('F41').Value = 120
$ws.Range('H41').Value = 'https://show.bilibili.com/platform/detail.html?id=86711'
$ws.Range('I41').Value = '//i2.hdslb.com/bfs/openplatform/202405/bllJHQFL1716983812432.jpeg'
$ws.Range('C42').Value = '上海·第十三届Redamancy动漫游戏嘉年华'
$ws.Range('D42').Value = '中山北路3300号 上海JOYPOLIS世嘉都市乐园'
$ws.Range('F42').Value = 2
$ws.Range('G42').Value = 60
$ws.Range('H42').Value = 'https://show.bilibili.com/platform/detail.html?id=87561'
$ws.Range('I42').Value = '//i2.hdslb.com/bfs/openplatform/202406/MkDJL3tI1718292674360.png'
$ws.Range('B43').Value = '2024-08-03'
$ws.Range('C43').Value = '上海·首届ICG动漫游戏博览会'
$ws.Range('D43').Value = '长江路258号C6 星球影棚'
$ws.Range('E43').Value = '2024.08.03 10:00-08.04 17:00'
$ws.Range('F43').Value = 2067
$ws.Range('G43').Value = 68
$ws.Range('H43').Value = 'https://show.bilibili.com/platform/detail.html?id=87419'
$ws.Range('I43').Value = '//i2.hdslb.com/bfs/openplatform/202406/thbceOGd1718290731704.jpeg'
$ws.Range('B44').Value = '2024-08-04'
$ws.Range('C44').Value = '上海·赛马娘ONLY·星夜天航'
$ws.Range('D44').Value = '漕宝路1688号 诺宝中心酒店'
$ws.Range('E44').Value = '2024.08.04 10:00-08.04 17:00'
$ws.Range('F44').Value = 93
$ws.Range('G44').Value = 80
$ws.Range('H44').Value = 'https://show.bilibili.com/platform/detail.html?id=87117'
$ws.Range('I44').Value = '//i2.hdslb.com/bfs/openplatform/202406/bdPnwqoq1717659799773.jpeg'
$ws.Range('C45').Value = '上海·创世次元动漫游戏嘉年华8.0'
$ws.Range('D45').Value = '漕宝路3366号 七宝万科广场'
$ws.Range('E45').Value = '2024.08.10 10:00-08.10 17:00'
$ws.Range('F45').Value = 839
$ws.Range('G45').Value = 58
$ws.Range('H45').Value = 'https://show.bilibili.com/platform/detail.html?id=86337'
$ws.Range('I45').Value = '//i1.hdslb.com/bfs/openplatform/202405/MryCBK6p1716736927373.jpeg'
$ws.Range('B46').Value = '2024-08-10'
$ws.Range('C46').Value = '上海·创造力动漫游戏嘉年华1.0'
$ws.Range('D46').Value = '莘福路288号 美莘商业广场'
$ws.Range('E46').Value = '2024.08.10 10:00-08.11 17:00'
$ws.Range('F46').Value = 8
$ws.Range('G46').Value = 65
$ws.Range('H46').Value = 'https://show.bilibili.com/platform/detail.html?id=87667'
$ws.Range('I46').Value = '//i2.hdslb.com/bfs/openplatform/202406/cjmOiK0E1718378936182.png'
$ws.Range('C47').Value = '上海·创世次元动漫游戏嘉年华9.0'
$ws.Range('D47').Value = '漕宝路3366号 七宝万科广场'
$ws.Range('E47').Value = '2024.08.17 10:00-08.17 17:00'
$ws.Range('F47').Value = 808
$ws.Range('G47').Value = 58
$ws.Range('H47').Value = 'https://show.bilibili.com/platform/detail.html?id=86355'
$ws.Range('I47').Value = '//i1.hdslb.com/bfs/openplatform/202405/UCMNMGbH1716782429997.jpeg'
$ws.Range('C48').Value = '上海·第六届燃梦BACG PRO动漫嘉年华-我们在燃梦相遇吧！'
$ws.Range('D48').Value = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws.Range('E48').Value = '2024.08.17 11:00-08.18 16:00'
$ws.Range('F48').Value = 121
$ws.Range('G48').Value = 65.8
$ws.Range('H48').Value = 'https://show.bilibili.com/platform/detail.html?id=85239'
$ws.Range('I48').Value = '//i1.hdslb.com/bfs/openplatform/202405/mzD4rhY21715109458100.jpeg'

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
$ws.Range('G5').Value = '不可售'
$ws.Range('F9').Value = 2606
$ws.Range('F10').Value = 1224
$ws.Range('F13').Value = 261
$ws.Range('F14').Value = 40
$ws.Range('F20').Value = 25
$ws.Range('F23').Value = 0
$ws.Range('F27').Value = 34
$ws.Range('F31').Value = 228
$ws.Range('F34').Value = 59
$ws.Range('F35').Value = 27
$ws.Range('F40').Value = 36
$ws.Range('F42').Value = 65
$ws.Range('F45').Value = 68

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 270
$ws.Range('F5').Value = 2922
$ws.Range('F6').Value = 4689
$ws.Range('F7').Value = 146
$ws.Range('F9').Value = 594
$ws.Range('F10').Value = 772
$ws.Range('F11').Value = 482
$ws.Range('F12').Value = 412
$ws.Range('F13').Value = 1147
$ws.Range('F14').Value = 316
$ws.Range('F15').Value = 752

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 1889
$ws.Range('F3').Value = 270
$ws.Range('F5').Value = 4689
$ws.Range('F6').Value = 772
$ws.Range('F8').Value = 412
$ws.Range('F9').Value = 412
$ws.Range('F11').Value = 928
$ws.Range('F12').Value = 1224
$ws.Range('F13').Value = 1626
$ws.Range('F14').Value = 1288
$ws.Range('F15').Value = 1559
$ws.Range('F16').Value = 71
$ws.Range('F17').Value = 1557
$ws.Range('F18').Value = 261
$ws.Range('F20').Value = 1701
$ws.Range('F21').Value = 1132
$ws.Range('F22').Value = 378
$ws.Range('F24').Value = 752
$ws.Range('F25').Value = 752
$ws.Range('F26').Value = 1813
$ws.Range('F27').Value = 252
$ws.Range('F28').Value = 825
$ws.Range('F29').Value = 8
$ws.Range('F31').Value = 1260
$ws.Range('F33').Value = 1077
$ws.Range('F34').Value = 82
$ws.Range('F35').Value = 1175
$ws.Range('F36').Value = 912
$ws.Range('F37').Value = 1178
$ws.Range('F39').Value = 1120
$ws.Range('F40').Value = 286
$ws.Range('F41').Value = 891
$ws.Range('F43').Value = 1699
$ws.Range('F46').Value = 2067
$ws.Range('F47').Value = 93
$ws.Range('F48').Value = 839
$ws.Range('F50').Value = 121
$ws.Range('F52').Value = 68
